$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 25.46092070582212
$ws.Range("C3").Value = 35.47070001895196
$ws.Range("C4").Value = 62.65030358103355
$ws.Range("C5").Value = 53.16831676755705
$ws.Range("C6").Value = 56.72933313736404
$ws.Range("C7").Value = 79.56627411049811
$ws.Range("C8").Value = 100
$ws.Range("C9").Value = 61.77368338896039
$ws.Range("C10").Value = 62.29318466346002
$ws.Range("C11").Value = 48.36269504309307
$ws.Range("C12").Value = 26.21767906357819
$ws.Range("C13").Value = 100
$ws.Range("C14").Value = 90.68989183504158
$ws.Range("C15").Value = 30.81115020868923
$ws.Range("C16").Value = 55.44047744661044
$ws.Range("C17").Value = 56.56134163456932
$ws.Range("C18").Value = 65.53142484668035
$ws.Range("C19").Value = 38.10080483931941
$ws.Range("C20").Value = 28.00689350224013
$ws.Range("C21").Value = 95.04943788441064
$ws.Range("C22").Value = 48.89928110391109
$ws.Range("C23").Value = 13.43946833874935
$ws.Range("C24").Value = 20.328262771856
$ws.Range("C25").Value = 15.38001080534552
$ws.Range("C26").Value = 80.46100733464016
$ws.Range("C27").Value = 19.48720233820424
$ws.Range("C28").Value = 14.68810980988851
$ws.Range("C29").Value = 72.70511039761431
$ws.Range("C30").Value = 51.51598588623294
$ws.Range("C31").Value = 24.49163739036445
$ws.Range("C32").Value = 43.97070949561198
$ws.Range("C33").Value = 76.01442758840712
$ws.Range("C34").Value = 58.30508627637893
$ws.Range("C35").Value = 48.02463459268451
$ws.Range("C36").Value = 41.027841706961
$ws.Range("C37").Value = 30.29428735738359
$ws.Range("C38").Value = 100
$ws.Range("C39").Value = 8.560731798254054
$ws.Range("C40").Value = 62.2724052047541
$ws.Range("C41").Value = 68.77969285782628
$ws.Range("C42").Value = 29.36823407648259
$ws.Range("C43").Value = 26.33800329683886
$ws.Range("C44").Value = 100
$ws.Range("C45").Value = 62.04430198085566
$ws.Range("C46").Value = 66.29812480115532
$ws.Range("C47").Value = 94.7708932802425
$ws.Range("C48").Value = 75.71909102655357
$ws.Range("C49").Value = 100
$ws.Range("C50").Value = 90.50399979915296
$ws.Range("C51").Value = 39.69592121582553
$ws.Range("C52").Value = 88.21132617235055
$ws.Range("C53").Value = 76.24009844570018
$ws.Range("C54").Value = 100
$ws.Range("C55").Value = 52.44895385076189
$ws.Range("C56").Value = 100
$ws.Range("C57").Value = 100
$ws.Range("C58").Value = 100
$ws.Range("C59").Value = 99.20256804176583
$ws.Range("C60").Value = 78.96361241406994
$ws.Range("C61").Value = 30.38223535434911
$ws.Range("C62").Value = 100
$ws.Range("C63").Value = 44.82603157826325
$ws.Range("C64").Value = 45.0693541582294
$ws.Range("C65").Value = 43.87679514285356
$ws.Range("C66").Value = 100
$ws.Range("C67").Value = 41.70673899721343
$ws.Range("C68").Value = 79.27647896967049
$ws.Range("C69").Value = 76.36980188576024
$ws.Range("C70").Value = 45.94876242539156
$ws.Range("C71").Value = 100
$ws.Range("C72").Value = 83.84160136738717
$ws.Range("C73").Value = 37.17206092743768
$ws.Range("C74").Value = 86.19931411257902
$ws.Range("C75").Value = 100
$ws.Range("C76").Value = 92.67614514625859
$ws.Range("C77").Value = 15.46691612318783
$ws.Range("C78").Value = 99.80490147995826
$ws.Range("C79").Value = 55.34113756524128
$ws.Range("C80").Value = 100
$ws.Range("C81").Value = 48.66982999772584
$ws.Range("C82").Value = 79.76887265481093
$ws.Range("C83").Value = 76.34167738253856
$ws.Range("C84").Value = 100
$ws.Range("C85").Value = 57.66466763358127
$ws.Range("C86").Value = 29.89297777248765
$ws.Range("C87").Value = 76.58343204936317
$ws.Range("C88").Value = 100
$ws.Range("C89").Value = 43.46939083426755
$ws.Range("C90").Value = 100
$ws.Range("C91").Value = 73.72837648363699
$ws.Range("C92").Value = 78.14059510703973
$ws.Range("C93").Value = 100
$ws.Range("C94").Value = 64.38484309842855
$ws.Range("C95").Value = 95.65773559565649
$ws.Range("C96").Value = 100
$ws.Range("C97").Value = 92.36382799119454
$ws.Range("C98").Value = 81.5009264733345
$ws.Range("C99").Value = 40.00522405864503
$ws.Range("C100").Value = 96.0828309977109
$ws.Range("C101").Value = 38.51221430445766
